$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd worksheet in this workbook
$ws3 = $wb.Worksheets.Item(3)

# Insert a new (empty) column before column N ("Late"), which pushes the
# existing N ("Late"), O (blank spacer) and P ("Outstanding") columns one
# position to the right (-> O, P, Q respectively).
$ws3.Columns("N").Insert()

# Re-select the sheet / cell that matches the post-edit selection
$ws3.Activate()
$ws3.Range("N9").Select()
